$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.775.28"
$ws.Range("E2").Value = "  -0.86%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.556.43"
$ws.Range("E3").Value = "  -1.54%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.09"
$ws.Range("E5").Value = "  -2.78%  "

$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "196.48"
$ws.Range("E6").Value = "  +0.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.612"
$ws.Range("E7").Value = "  -2.28%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.209"
$ws.Range("E9").Value = "  +0.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.627"
$ws.Range("E10").Value = "  -3.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.76"
$ws.Range("E11").Value = "  -1.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000290"
$ws.Range("E12").Value = "  -4.61%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.26"
$ws.Range("E13").Value = "  -3.19%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.118.48"
$ws.Range("E14").Value = "  -1.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "663.94"
$ws.Range("E15").Value = "  +11.32%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.830.37"
$ws.Range("E16").Value = "  -0.85%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "12.57"
$ws.Range("E17").Value = "  -3.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.548.45"
$ws.Range("E18").Value = "  -1.78%  "

$ws.Range("E19").Value = "  -0.99%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.46"
$ws.Range("E20").Value = "  -3.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.965"
$ws.Range("E21").Value = "  -3.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.39"
$ws.Range("E22").Value = "  +2.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.37"
$ws.Range("E23").Value = "  +3.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "105.45"
$ws.Range("E24").Value = "  +3.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.39"
$ws.Range("E25").Value = "  -5.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.93"
$ws.Range("E26").Value = "  -3.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.20"
$ws.Range("E27").Value = "  -5.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.63"
$ws.Range("E28").Value = "  +0.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.40"
$ws.Range("E29").Value = "  -1.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.45"
$ws.Range("E30").Value = "  -6.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.82"
$ws.Range("E31").Value = "  -5.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.79"
$ws.Range("E32").Value = "  -4.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.112"
$ws.Range("E33").Value = "  -4.92%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "61.76"
$ws.Range("E34").Value = "  -2.67%  "

$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.761.07"
$ws.Range("E35").Value = "  -3.67%  "

$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0826"
$ws.Range("E36").Value = "  -8.14%  "

$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.77"
$ws.Range("E37").Value = "  +6.76%  "

$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "515.23"
$ws.Range("E39").Value = "  -5.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.94"
$ws.Range("E40").Value = "  -6.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.373"
$ws.Range("E41").Value = "  -4.62%  "

$ws.Range("E42").Value = "  +0.76%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "34.75"
$ws.Range("E43").Value = "  -5.94%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0456"
$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("E45").Value = "  -0.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.89"
$ws.Range("E46").Value = "  +0.72%  "

$ws.Range("E47").Value = "  -2.17%  "

$ws.Range("E48").Value = "  -0.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.36"
$ws.Range("E49").Value = "  -2.83%  "

$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000239"
$ws.Range("E50").Value = "  -5.03%  "

$ws.Range("E51").Value = "  +18.45%  "
